# Update countries & provincias Spain
# Applies the data refresh described by the commit:
#  - Nigeria and Liberia move up in the (descending, by total cases) ranking,
#    so the rows that used to hold Armenia/Irak/Nigeria and Cabo Verde/Liberia
#    now hold Nigeria/Armenia/Irak and Liberia/Cabo Verde respectively.
#  - A handful of other numeric updates (Estados Unidos, Canada).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4): refreshed totals ---
$ws.Cells.Item(4, 2).Value = 1159430   # Casos totales
$ws.Cells.Item(4, 3).Value = 28400     # Nuevos casos
$ws.Cells.Item(4, 5).Value = 931371    # Recuperados
$ws.Cells.Item(4, 7).Value = 1638      # Casos criticos
$ws.Cells.Item(4, 8).Value = 67391     # Muertes

# --- Canada (row 15): refreshed totals ---
$ws.Cells.Item(15, 4).Value = 23801    # Casos activos
$ws.Cells.Item(15, 5).Value = 29347    # Recuperados
$ws.Cells.Item(15, 7).Value = 175      # Casos criticos
$ws.Cells.Item(15, 8).Value = 3566     # Muertes

# --- Nigeria climbs above Armenia and Irak (rows 68-70) ---
$ws.Cells.Item(68, 1).Value = "Nigeria"
$ws.Cells.Item(68, 2).Value = 2388
$ws.Cells.Item(68, 3).Value = 218
$ws.Cells.Item(68, 4).Value = 351
$ws.Cells.Item(68, 5).Value = 1952
$ws.Cells.Item(68, 6).Value = 2
$ws.Cells.Item(68, 7).Value = 17
$ws.Cells.Item(68, 8).Value = 85

$ws.Cells.Item(69, 1).Value = "Armenia"
$ws.Cells.Item(69, 2).Value = 2273
$ws.Cells.Item(69, 3).Value = 125
$ws.Cells.Item(69, 4).Value = 1010
$ws.Cells.Item(69, 5).Value = 1230
$ws.Cells.Item(69, 6).Value = 10
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 33

$ws.Cells.Item(70, 1).Value = "Irak"
$ws.Cells.Item(70, 2).Value = 2219
$ws.Cells.Item(70, 3).Value = 66
$ws.Cells.Item(70, 4).Value = 1473
$ws.Cells.Item(70, 5).Value = 651
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 95

# --- Liberia climbs above Cabo Verde (rows 138-139) ---
$ws.Cells.Item(138, 1).Value = "Liberia"
$ws.Cells.Item(138, 2).Value = 154
$ws.Cells.Item(138, 3).Value = 2
$ws.Cells.Item(138, 4).Value = 48
$ws.Cells.Item(138, 5).Value = 88
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 18

$ws.Cells.Item(139, 1).Value = "Cabo Verde"
$ws.Cells.Item(139, 2).Value = 152
$ws.Cells.Item(139, 3).Value = 30
$ws.Cells.Item(139, 4).Value = 18
$ws.Cells.Item(139, 5).Value = 132
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = 2
